# Update cryptocurrency price/volume figures in the worksheet.
# Source values are stored as text (t="inlineStr"/shared-string) even when
# they look numeric (e.g. "0.3667", "20.559.70"), so every assignment is
# apostrophe-prefixed to force Excel to keep it as text instead of silently
# parsing it into a floating point number (which would also mangle things
# like trailing zeros -> "0.00001030" becoming "1.03E-05"). The Style reset
# afterwards clears the quote-prefix flag so the cell ends up with the same
# default "Normal" styling as before the edit, without a literal "'" showing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'20.559.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.22%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.480.24"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.62%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.07%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("E5").Value = "'  +2.48%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'279.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.96%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.3667"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -1.56%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3081"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -3.72%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'39.95"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -4.69%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -0.58%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.06663"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.87%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +0.03%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'5.510"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -2.37%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'18.11"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.58%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'6.207"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.38%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.9780"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +2.40%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.00001030"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.17%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'1.481.59"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.33%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.05940"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +2.81%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'69.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -4.71%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'5.480"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -4.31%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -2.41%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'11.04"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.61%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.253"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.78%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'20.622.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.53%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'141.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +2.43%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -7.61%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -2.26%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.640.20"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.04%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'113.68"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.37%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'3.958"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.42%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'5.008"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -6.87%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.8186"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -2.97%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.08029"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +1.90%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.549"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -6.34%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'1.231"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +9.40%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.05820"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -4.89%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'4.730"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -4.42%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.9775"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +1.41%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'7.729"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +6.29%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -1.58%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -3.19%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.1889"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.20%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.5299"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -2.92%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -1.92%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'12.24"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -2.49%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'119.27"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -2.06%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.5209"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -3.32%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.802"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -2.00%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.06471"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.16%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.9916"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.21%  "
$ws.Range("E51").Style = "Normal"
